$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update existing "Test2" block (rows 2 and 5) with new dynamic values ---
# Row 2 (first block's value row)
$ws.Range("C2").Value = "'749939"
$ws.Range("E2").Value = "'749939"

# Row 5 (second block's value row, the one that gets duplicated below)
$ws.Range("C5").Value = "'749939"
$ws.Range("E5").Value = "'749939"
$ws.Range("F5").Value = "'749933"
$ws.Range("H5").Value = "'749933"
$ws.Range("I5").Value = "'749936"
$ws.Range("K5").Value = "'749936"

# --- Duplicate the "Test2" block (rows 4:6) into a new "Test3" block (rows 7:9) ---
$ws.Range("A4:L6").Copy($ws.Range("A7"))

# The paste fills the whole destination rectangle; trim back the cells that
# must stay entirely empty so the row layout matches rows 4:6 exactly.
$ws.Range("L7").ClearContents()
$ws.Range("A8").ClearContents()
$ws.Range("L8").ClearContents()
$ws.Range("A9").ClearContents()
$ws.Range("B9").ClearContents()
$ws.Range("D9").ClearContents()
$ws.Range("I9").ClearContents()
$ws.Range("J9").ClearContents()
$ws.Range("K9").ClearContents()

# Relabel the new block
$ws.Range("A7").Value = "Test3"
$ws.Range("L9").Value = "Test3"

# New hyperlink for the duplicated URL cell (Hyperlinks.Add resets the cell
# style, so restore the same "Hyperlink" formatting used by B2/B5 afterwards)
$ws.Hyperlinks.Add($ws.Range("B8"), "https://localhost:8080/")
$ws.Range("B5").Copy()
$ws.Range("B8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Column width tweaks ---
$ws.Columns.Item(3).ColumnWidth = 13.166666666666666
$ws.Columns.Item(6).ColumnWidth = 11.416666666666666

# --- Selection bookkeeping ---
[void]$ws.Range("E12").Select()

$wb.Save()
